# Commit: [main] Set prescaler and counter period of TIM21(1ms) and TIM22(1us) below.
# - TIM21, 1ms, prescaler[3199], counter period[9]
# - TIM22, 1us, prescaler[31], counter period[65535]
# - Verify by __MEASURE_TIM22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWM")

# The old "measure execution time" table (D23:G32) is removed entirely - it's
# replaced by two new repeated Clock/Prescaler/fPWM/Counter-Period blocks
# (TIM21 @ 1ms and TIM22 @ 1us) below the existing base block.
$ws.Range("D23:G32").Clear()

# The base block's 1/fPWM "period" helper formulas (D12, D13:E13) move down
# to the new 1us block (D24, D25:E25) instead - clear them here.
$ws.Range("D12:E13").Clear()

# Clone the formatting of the existing base PWM-setting block (rows 10-13,
# styles s=2/s=3 with fill+border) down into the two new blocks so every
# block looks the same.
$ws.Range("B10:C13").Copy()
$ws.Range("B16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B10:C13").Copy()
$ws.Range("B22").PasteSpecial(-4122)   # xlPasteFormats

# Clone the section-header formatting (s=1, bold, no fill/border) for the new
# "1ms" / "1us" sub-headers.
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B9").Copy()
$ws.Range("B21").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Base PWM setting block (rows 10-13) -------------------------------
$ws.Range("B10").Value = "Clock"
$ws.Range("C10").Formula = "=32*1000*1000"

$ws.Range("B11").Value = "Prescaler"
$ws.Range("C11").Value = 0

$ws.Range("C12").Formula = "=C10/(C13+1)/(C11+1)"

$ws.Range("B13").Value = "Counter Period"
$ws.Range("C13").Value = 31999

# --- TIM21, 1ms block (rows 15-19) --------------------------------------
$ws.Range("B15").Value = "1ms"

$ws.Range("B16").Value = "Clock"
$ws.Range("C16").Formula = "=32*1000*1000"

$ws.Range("B17").Value = "Prescaler"
$ws.Range("C17").Value = 3199

$ws.Range("C18").Formula = "=C16/(C19+1)/(C17+1)"

$ws.Range("B19").Value = "Counter Period"
$ws.Range("C19").Value = 9

# --- TIM22, 1us block (rows 21-25) --------------------------------------
$ws.Range("B21").Value = "1us"

$ws.Range("B22").Value = "Clock"
$ws.Range("C22").Formula = "=32*1000*1000"

$ws.Range("B23").Value = "Prescaler"
$ws.Range("C23").Value = 31

$ws.Range("C24").Formula = "=C22/(C25+1)/(C23+1)"
$ws.Range("D24").Formula = "=1/C24"

$ws.Range("B25").Value = "Counter Period"
$ws.Range("C25").Value = 65535
$ws.Range("D25").Formula = "=D24/C25"
$ws.Range("E25").Formula = "=D25*1000"

# --- fPWM(Hz) labels - rich text: "f" + subscript-bold "PWM" + bold "(Hz)" --
foreach ($addr in @("B12", "B18", "B24")) {
    $rng = $ws.Range($addr)
    $rng.Value = "fPWM(Hz)"
    $c = $rng.Characters(1, 1)
    $c.Font.Bold = $true
    $c = $rng.Characters(2, 3)
    $c.Font.Bold = $true
    $c.Font.Subscript = $true
    $c = $rng.Characters(5, 4)
    $c.Font.Bold = $true
}

$ws.Rows.Item(12).RowHeight = 18
$ws.Rows.Item(18).RowHeight = 18
$ws.Rows.Item(24).RowHeight = 18

# --- Selection / view state ---------------------------------------------
$ws.Activate()
$ws.Range("E28").Select()
